# Apply the edit described by the diff:
# - In the second results table (rows 11-17, "Section 7 results"), rename the
#   "Validation Score" column headers (C13 and E13) to "Test Score".
# - Move the active cell selection to H13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C13").Value = "Test Score"
$ws.Range("E13").Value = "Test Score"

$ws.Range("H13").Select()
